# Update column F ("dSF") values on a handful of rows to reflect the
# repulled / recalculated data (commit: "repull data, push all data, mean calculation").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -3
$ws.Range("F3").Value  = -8
$ws.Range("F8").Value  = -3
$ws.Range("F9").Value  = 1
$ws.Range("F11").Value = 0
$ws.Range("F17").Value = 4
$ws.Range("F21").Value = -6
$ws.Range("F22").Value = -12
$ws.Range("F26").Value = -5
$ws.Range("F28").Value = -4
$ws.Range("F32").Value = 1
$ws.Range("F33").Value = -11
$ws.Range("F34").Value = 1
